$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.586.01'
$ws.Range('E2').Value = '  -1.34%  '
$ws.Range('D3').Value = '2.285.85'
$ws.Range('E3').Value = '  +1.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '95.05'
$ws.Range('E5').Value = '  -3.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '266.72'
$ws.Range('E6').Value = '  -2.63%  '
$ws.Range('E7').Value = '  -1.06%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -3.83%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '44.58'
$ws.Range('E10').Value = '  -7.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0934'
$ws.Range('E11').Value = '  -1.25%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.74'
$ws.Range('E12').Value = '  -5.41%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.105'
$ws.Range('E13').Value = '  -0.07%  '
$ws.Range('D14').Value = '2.627.05'
$ws.Range('E14').Value = '  +1.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.14'
$ws.Range('E15').Value = '  -2.48%  '
$ws.Range('E16').Value = '  +0.61%  '
$ws.Range('D17').Value = '2.285.94'
$ws.Range('E17').Value = '  +1.03%  '
$ws.Range('D18').Value = '43.508.57'
$ws.Range('E18').Value = '  -1.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000107'
$ws.Range('E19').Value = '  -0.49%  '
$ws.Range('E20').Value = '  -1.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.32'
$ws.Range('E21').Value = '  +1.85%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.42'
$ws.Range('E22').Value = '  +2.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '234.40'
$ws.Range('E23').Value = '  -0.39%  '
$ws.Range('E24').Value = '  -15.36%  '
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('E26').Value = '  -1.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.17'
$ws.Range('E27').Value = '  -2.49%  '
$ws.Range('E28').Value = '  -0.68%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '39.88'
$ws.Range('E29').Value = '  -1.27%  '
$ws.Range('E30').Value = '  -0.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '175.67'
$ws.Range('E31').Value = '  +1.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.86'
$ws.Range('E32').Value = '  +3.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0880'
$ws.Range('E33').Value = '  -4.36%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.33'
$ws.Range('E34').Value = '  -7.29%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  -6.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0355'
$ws.Range('E37').Value = '  -0.78%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.41'
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.31'
$ws.Range('E39').Value = '  -6.87%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.34'
$ws.Range('E40').Value = '  +7.15%  '
$ws.Range('E41').Value = '  -7.64%  '
$ws.Range('E42').Value = '  +15.57%  '
$ws.Range('B43').Value = 'MultiversX'
$ws.Range('C43').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '63.93'
$ws.Range('E43').Value = '  +2.14%  '
$ws.Range('B44').Value = 'Celestia'
$ws.Range('C44').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '11.94'
$ws.Range('E44').Value = '  -5.58%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.80'
$ws.Range('E45').Value = '  +2.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.22'
$ws.Range('E46').Value = '  -4.83%  '
$ws.Range('E47').Value = '  -1.99%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '97.76'
$ws.Range('E48').Value = '  -2.91%  '
$ws.Range('E49').Value = '  -1.11%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.506.82'
$ws.Range('E50').Value = '  +1.16%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.49'
$ws.Range('E51').Value = '  +4.24%  '
